$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 7 de Agosto de 2020 a las 22:58"

# Row 4: Estados Unidos
$ws.Range("B4").Value = 5077556
$ws.Range("C4").Value = 45377
$ws.Range("D4").Value = 2591656
$ws.Range("E4").Value = 2322254
$ws.Range("G4").Value = 842
$ws.Range("H4").Value = 163646

# Row 8: Sudafrica
$ws.Range("B8").Value = 545476
$ws.Range("C8").Value = 7292
$ws.Range("D8").Value = 394759
$ws.Range("E8").Value = 140808
$ws.Range("G8").Value = 305
$ws.Range("H8").Value = 9909

# Row 12: Colombia -> España
$ws.Range("A12").Value = "España"
$ws.Range("B12").Value = 361442
$ws.Range("C12").Value = 4507
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 0
$ws.Range("G12").Value = 3
$ws.Range("H12").Value = 28503

# Row 13: España -> Colombia
$ws.Range("A13").Value = "Colombia"
$ws.Range("B13").Value = 357710
$ws.Range("D13").Value = 192355
$ws.Range("E13").Value = 153416
$ws.Range("H13").Value = 11939

# Row 22: Alemania
$ws.Range("B22").Value = 216315
$ws.Range("C22").Value = 1105
$ws.Range("E22").Value = 10861
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 9254

# Row 23: Francia
$ws.Range("D23").Value = 82836
$ws.Range("E23").Value = 84761

# Row 52: Barein
$ws.Range("E52").Value = 2785
$ws.Range("G52").Value = 3
$ws.Range("H52").Value = 159

# Row 59: Argelia
$ws.Range("B59").Value = 34155
$ws.Range("C59").Value = 529
$ws.Range("D59").Value = 23667
$ws.Range("E59").Value = 9206
$ws.Range("G59").Value = 9
$ws.Range("H59").Value = 1282

# Row 66: Kenia
$ws.Range("B66").Value = 25138
$ws.Range("C66").Value = 727
$ws.Range("D66").Value = 11118
$ws.Range("E66").Value = 13607
$ws.Range("G66").Value = 14
$ws.Range("H66").Value = 413

# Row 90: Gabon
$ws.Range("B90").Value = 7923
$ws.Range("C90").Value = 136
$ws.Range("D90").Value = 5704
$ws.Range("E90").Value = 2168

# Row 91: Tayikistan -> Guinea
$ws.Range("A91").Value = "Guinea"
$ws.Range("B91").Value = 7777
$ws.Range("C91").Value = 113
$ws.Range("D91").Value = 6800
$ws.Range("E91").Value = 927
$ws.Range("G91").Value = 1
$ws.Range("H91").Value = 50

# Row 92: Guinea -> Tayikistan
$ws.Range("A92").Value = "Tayikistan"
$ws.Range("B92").Value = 7706
$ws.Range("C92").Value = 41
$ws.Range("D92").Value = 6484
$ws.Range("E92").Value = 1160
$ws.Range("H92").Value = 62

# Row 106: Maldivas
$ws.Range("B106").Value = 4769
$ws.Range("C106").Value = 89
$ws.Range("D106").Value = 2754
$ws.Range("E106").Value = 1996

# Row 107: Hungria -> Republica de Africa Central
$ws.Range("A107").Value = "Republica de Africa Central"
$ws.Range("B107").Value = 4622
$ws.Range("C107").Value = 2
$ws.Range("D107").Value = 1672
$ws.Range("E107").Value = 2891
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = 59

# Row 108: Republica de Africa Central -> Hungria
$ws.Range("A108").Value = "Hungria"
$ws.Range("B108").Value = 4621
$ws.Range("C108").Value = 24
$ws.Range("D108").Value = 3464
$ws.Range("E108").Value = 555
$ws.Range("G108").Value = 2
$ws.Range("H108").Value = 602

# Row 151: Togo
$ws.Range("B151").Value = 1028
$ws.Range("C151").Value = 16
$ws.Range("D151").Value = 710
$ws.Range("E151").Value = 296

# Row 202: Santa Lucia -> Timor Oriental
$ws.Range("A202").Value = "Timor Oriental"

# Row 203: Timor Oriental -> Santa Lucia
$ws.Range("A203").Value = "Santa Lucia"
